$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 954.7912
$ws.Range("I15").Value = 954.7912
$ws.Range("K15").Value = 2864.3736
$ws.Range("M15").Value = -2695.3736
# Row 19
$ws.Range("H19").Value = 3760857.2
$ws.Range("I19").Value = 6579944.5
$ws.Range("J19").Value = 2074
$ws.Range("K19").Value = 6579944.5
$ws.Range("L19").Value = 2074
$ws.Range("M19").Value = -6579769.5
$ws.Range("N19").Value = -2424
# Row 76
$ws.Range("H76").Value = 3327.2727
$ws.Range("I76").Value = 3277.7778
$ws.Range("K76").Value = 3277.7778
$ws.Range("M76").Value = -2962.7778
# Row 79
$ws.Range("H79").Value = 3327.2727
$ws.Range("I79").Value = 3277.7778
$ws.Range("K79").Value = 3277.7778
$ws.Range("M79").Value = -2185.7778
# Row 129
$ws.Range("H129").Value = 821
$ws.Range("I129").Value = 342.30768
$ws.Range("J129").Value = 972.7805
$ws.Range("K129").Value = 1026.92304
$ws.Range("L129").Value = 2918.3415
$ws.Range("M129").Value = 3973.07696
$ws.Range("N129").Value = -12918.3415
# Row 132
$ws.Range("H132").Value = 21493396
$ws.Range("I132").Value = 24394016
$ws.Range("K132").Value = 73182048
$ws.Range("M132").Value = -73179518
# Row 135
$ws.Range("H135").Value = 409.76666
$ws.Range("I135").Value = 384.3158
$ws.Range("J135").Value = 453.72726
$ws.Range("K135").Value = 3458.8422
$ws.Range("L135").Value = 4083.54534
$ws.Range("M135").Value = -923.8422
$ws.Range("N135").Value = -9153.545340000001
# Row 137
$ws.Range("H137").Value = 3994.889
$ws.Range("I137").Value = 2586.2856
$ws.Range("J137").Value = 4487.9
$ws.Range("K137").Value = 7758.8568
$ws.Range("L137").Value = 13463.7
$ws.Range("M137").Value = -5208.8568
$ws.Range("N137").Value = -18563.7
# Row 138
$ws.Range("H138").Value = 4675.44
$ws.Range("I138").Value = 943.8261
$ws.Range("J138").Value = 5790.078
$ws.Range("K138").Value = 2831.4783
$ws.Range("L138").Value = 17370.234
$ws.Range("M138").Value = 2308.5217
$ws.Range("N138").Value = -27650.234
# Row 141
$ws.Range("H141").Value = 26266.65
$ws.Range("I141").Value = 27422.79
$ws.Range("J141").Value = 4300
$ws.Range("K141").Value = 82268.37
$ws.Range("L141").Value = 12900
$ws.Range("M141").Value = -77088.37
$ws.Range("N141").Value = -23260

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 499.5357
$ws.Range("I2").Value = 540
$ws.Range("J2").Value = 398.375
$ws.Range("K2").Value = 540
$ws.Range("L2").Value = 398.375
$ws.Range("M2").Value = -427
$ws.Range("N2").Value = -624.375
# Row 32
$ws.Range("H32").Value = 4403.7456
$ws.Range("I32").Value = 4121.0815
$ws.Range("K32").Value = 4121.0815
$ws.Range("M32").Value = -3834.0815
# Row 45
$ws.Range("H45").Value = 1050
$ws.Range("I45").Value = 1050
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1050
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -673
$ws.Range("N45").ClearContents()
# Row 74
$ws.Range("H74").Value = 3494.7568
$ws.Range("I74").Value = 3469.9333
$ws.Range("K74").Value = 3469.9333
$ws.Range("M74").Value = -2595.9333
# Row 77
$ws.Range("H77").Value = 3494.7568
$ws.Range("I77").Value = 3469.9333
$ws.Range("K77").Value = 17349.6665
$ws.Range("M77").Value = -12981.6665
# Row 116
$ws.Range("H116").Value = 499.5357
$ws.Range("I116").Value = 540
$ws.Range("J116").Value = 398.375
$ws.Range("K116").Value = 540
$ws.Range("L116").Value = 398.375
$ws.Range("M116").Value = 1754
$ws.Range("N116").Value = -4986.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 499.5357
$ws.Range("I3").Value = 540
$ws.Range("J3").Value = 398.375
$ws.Range("K3").Value = 540
$ws.Range("L3").Value = 398.375
$ws.Range("M3").Value = -426
$ws.Range("N3").Value = -626.375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11366896
$ws.Range("I31").Value = 1854.2916
$ws.Range("K31").Value = 1854.2916
$ws.Range("M31").Value = -1559.2916
# Row 34
$ws.Range("H34").Value = 11366896
$ws.Range("I34").Value = 1854.2916
$ws.Range("K34").Value = 1854.2916
$ws.Range("M34").Value = -1652.2916
# Row 58
$ws.Range("H58").Value = 1581.6428
$ws.Range("I58").Value = 1494.2174
$ws.Range("J58").Value = 1789.6552
$ws.Range("K58").Value = 1494.2174
$ws.Range("L58").Value = 1789.6552
$ws.Range("M58").Value = -1291.2174
$ws.Range("N58").Value = -2195.6552
# Row 136
$ws.Range("H136").Value = 1581.6428
$ws.Range("I136").Value = 1494.2174
$ws.Range("J136").Value = 1789.6552
$ws.Range("K136").Value = 4482.6522
$ws.Range("L136").Value = 5368.9656
$ws.Range("M136").Value = -1932.6522
$ws.Range("N136").Value = -10468.9656

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 531
$ws.Range("J113").Value = 597.84375
$ws.Range("L113").Value = 1793.53125
$ws.Range("N113").Value = -6133.53125
# Row 118
$ws.Range("H118").Value = 462.1
$ws.Range("I118").Value = 462.1
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1386.3
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -143.3000000000002
$ws.Range("N118").ClearContents()
# Row 131
$ws.Range("H131").Value = 835.1177
$ws.Range("I131").Value = 500.625
$ws.Range("J131").Value = 897.3488
$ws.Range("K131").Value = 1501.875
$ws.Range("L131").Value = 2692.0464
$ws.Range("M131").Value = 3538.125
$ws.Range("N131").Value = -12772.0464

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6417.4736
$ws.Range("I70").Value = 5823.423
$ws.Range("J70").Value = 7704.5835
$ws.Range("K70").Value = 5823.423
$ws.Range("L70").Value = 7704.5835
$ws.Range("M70").Value = -5553.423
$ws.Range("N70").Value = -8244.583500000001
# Row 73
$ws.Range("H73").Value = 6417.4736
$ws.Range("I73").Value = 5823.423
$ws.Range("J73").Value = 7704.5835
$ws.Range("K73").Value = 5823.423
$ws.Range("L73").Value = 7704.5835
$ws.Range("M73").Value = -4887.423
$ws.Range("N73").Value = -9576.583500000001
# Row 97
$ws.Range("H97").Value = 915
$ws.Range("I97").Value = 915
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 915
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -419
$ws.Range("N97").ClearContents()
# Row 122
$ws.Range("H122").Value = 3035.2222
$ws.Range("I122").Value = 1576.5
$ws.Range("J122").Value = 5952.6665
$ws.Range("K122").Value = 4729.5
$ws.Range("L122").Value = 17857.9995
$ws.Range("M122").Value = -2279.5
$ws.Range("N122").Value = -22757.9995
# Row 132
$ws.Range("H132").Value = 2208.9622
$ws.Range("I132").Value = 1355.973
$ws.Range("K132").Value = 4067.919
$ws.Range("M132").Value = -1537.919

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 7033.75
$ws.Range("I122").Value = 2700
$ws.Range("K122").Value = 8100
$ws.Range("M122").Value = -5650
# Row 132
$ws.Range("H132").Value = 8951.48
$ws.Range("I132").Value = 9234.387000000001
$ws.Range("J132").Value = 8489.895
$ws.Range("K132").Value = 27703.161
$ws.Range("L132").Value = 25469.685
$ws.Range("M132").Value = -25173.161
$ws.Range("N132").Value = -30529.685
